# Sync non-localizable rule data (CodeQuality-rules-latest-CS.xlsx)
#
# Summary of change:
#  - The "BannedPaths" rule row is removed from its original position
#    (row 35), and the following index/Search-related rows shift up.
#  - A renamed "BannedPath" rule (same description, but Severity raised
#    from Blocker to Critical, and no Tags) is (re-)inserted lower in the
#    table, at what is now row 40.
#  - The active cell selection on the sheet moves from B38 to A37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "BannedPaths" row; everything below ripples up one row.
$ws.Rows(35).Delete()

# Make room for the renamed rule lower in the list (now at row 40) by
# inserting a fresh blank row there, pushing "AEM Rules:AEM-3" (and
# everything after) back down.
$ws.Rows(40).Insert()

$ws.Range("A40").Value = "BannedPath"
$ws.Range("B40").Value = "Customer packages should not install content under /libs"
$ws.Range("C40").Value = "Bug"
$ws.Range("D40").Value = "Critical"

# Update the sheet's saved selection/active cell.
$ws.Range("A37").Select() | Out-Null
